$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7177848219871521
$ws.Range("B1").Value = 3.934299945831299
$ws.Range("C1").Value = 5.767133235931396
$ws.Range("D1").Value = 1.617388725280762
$ws.Range("E1").Value = 0.9706059098243713
